# "combine forgot_password and sign_up links on the login page to one test
#  case -> test_links"
#
# On the test_cases sheet, rows 7 and 8 are the login-page test cases
# "test_sign_up_link" and "test_forgot_password_link". They get merged into
# a single row, "test_page_links" / "Check if all page links are working",
# and the table (autofilter, filter-database name, selection, dimension)
# shrinks by the one row that disappears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_cases")

# Turn row 7 (currently "test_sign_up_link") into the merged test case.
$ws.Range("C7").Value = "test_page_links"
$ws.Range("D7").Value = "Check if all page links are working"

# Row 8 ("test_forgot_password_link") is now redundant. Deleting it shifts
# rows 9-11 up by one, so the former row 11 (profile page / test_changing_profile)
# becomes the new row 10, and the sheet shrinks to A2:H10.
$ws.Rows.Item(8).Delete() | Out-Null

# Re-apply the autofilter over the new, smaller range (A2:H11 -> A2:H10).
$ws.AutoFilterMode = $false
$ws.Range("A2:H10").AutoFilter() | Out-Null

# Keep the workbook's hidden _xlnm._FilterDatabase name in sync with the
# shrunken filter range.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=test_cases!`$A`$2:`$H`$10"

# Restore the saved cursor/selection position recorded for the sheet.
$ws.Range("E12").Select() | Out-Null
